# Insert a new record row at row 43 (pushing the existing rows 43-121 down
# to 44-122) and populate it with the new observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(43).Insert()

$ws.Range("A43").Value = 11
$ws.Range("B43").Value = "Vega Monumental Concepción"
$ws.Range("C43").Value = "Bíobío"
$ws.Range("D43").Value = 44775
$ws.Range("E43").Value = 8
$ws.Range("F43").Value = 100112021
$ws.Range("G43").Value = "Ají"
$ws.Range("H43").Value = "Inferno"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 35
$ws.Range("K43").Value = 11000
$ws.Range("L43").Value = 12000
$ws.Range("M43").Value = 11571
$ws.Range("N43").Value = "$/caja 12 kilos"
$ws.Range("O43").Value = "Región de Arica y Parinacota"
$ws.Range("P43").Value = 964
$ws.Range("Q43").Value = 12
$ws.Range("R43").Value = "Hortaliza"
